$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Footprint fix for JLCPCB assembly: U1's footprint changes from SOP-16 to TSSOP-20
$ws.Range("C2").Value = "TSSOP-20"

# Keep selection consistent with the saved workbook state
$ws.Range("C3").Select()
